$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column D, shifting the existing data (D..K) to (E..L)
$ws.Columns.Item(4).Insert()

# Copy number formatting from column E (the old column D, now shifted right)
# into the new blank column D, bounded to the populated range so the sheets
# UsedRange / dimension stays correct (whole-column copy would blow it out).
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new column D with the latest (FY2018) financial data
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 2189100
$ws.Range("D9").Value = 1608300
$ws.Range("D10").Value = 580800
$ws.Range("D12").Value = 55900
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 1817900
$ws.Range("D18").Value = 371200
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = 494300
$ws.Range("D22").Value = 37700
$ws.Range("D23").Value = 333500
$ws.Range("D24").Value = 69900
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 263600
$ws.Range("D27").Value = 269200
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 7400
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = 276600
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 276600
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 32700
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 311400
$ws.Range("D44").Value = 297800
$ws.Range("D45").Value = 33900
$ws.Range("D46").Value = 675800
$ws.Range("D47").Value = 48700
$ws.Range("D48").Value = 1876500
$ws.Range("D49").Value = 142300
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 80800
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 2824100
$ws.Range("D57").Value = 161900
$ws.Range("D58").Value = 9400
$ws.Range("D59").Value = 155400
$ws.Range("D60").Value = 326700
$ws.Range("D61").Value = 947400
$ws.Range("D62").Value = 228000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1502100
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 1726500
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1322000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 276600
$ws.Range("D83").Value = 123100
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 421400
$ws.Range("D91").Value = -184100
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -187500
$ws.Range("D96").Value = -48400
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -257300
$ws.Range("D101").Value = -4000
$ws.Range("D102").Value = -27400
